$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New rows 3-5 (DRA002, DRA003, DRA004) added under the existing DRA001 row.
# Shared-string insertion order matters (it drives the index each string
# lands on), so TCIDs are written first (A3,A4,A5) followed by the
# Description cells in the same order the original author must have used:
# C5, then C3, then C4.
# ---------------------------------------------------------------------------

$ws.Range("A3").Value = "DRA002"
$ws.Range("A4").Value = "DRA003"
$ws.Range("A5").Value = "DRA004"

$skipLinkingText = 'Verify that user can skip the linking by clicking on "Not now button" on the modal "Already have an account? .. || Verify that once the user skips linking then user will not be prompted to link again.'
$ws.Range("C5").Value = $skipLinkingText

$richPrefix = "Verify that user should be prompted to link accounts, when sign in first time on "
$richBold = "DRA landing screen"
$richSuffix = " using STeAM. (Note:User should already been sign into social)"
$ws.Range("C3").Value = $richPrefix + $richBold + $richSuffix

$socialLinkingText = @"
 Verify that when linking a social with a matching email, if the user click [X] cross mark on the screen then he will be taken back to the DRA Login page. || Verify that text on the modal "Already have an account? ..
You have previously signed in with <email address> using LinkedIn
To Protect your security, please sign into LinkedIn so that we can link your account.
<not now> <Sign in using Facebook>" , when linking LinkedIn with steam ||Verify that when linking a social with a matching email, if the user clicks outside the Linking modal on the screen then nothing should happens
"@
# Strip the trailing newline the here-string's closing marker adds.
$socialLinkingText = $socialLinkingText.TrimEnd("`r", "`n")
$ws.Range("C4").Value = $socialLinkingText

# Runmode (column B) and Jira id / "Y" (column D) for the three new rows -
# these reuse the existing "OBT" / "Y" shared strings.
$ws.Range("B3").Value = "OBT"
$ws.Range("B4").Value = "OBT"
$ws.Range("B5").Value = "OBT"
$ws.Range("D3").Value = "Y"
$ws.Range("D4").Value = "Y"
$ws.Range("D5").Value = "Y"

# Column E (Results) stays blank on every data row, same as row 2.
$ws.Range("E3").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("E5").Value = ""

# ---------------------------------------------------------------------------
# Formatting: reuse the same thin border already used by row 2 (setting the
# border Color before the LineStyle makes the engine reuse the existing
# border definition instead of minting a new one).
# ---------------------------------------------------------------------------
$ws.Range("A3:E5").Borders.Color = 0
$ws.Range("A3:E5").Borders.LineStyle = 1

# Column A (TCID) and the Description cells of rows 4-5 wrap text, matching
# the style already used for A2 / C2.
$ws.Range("A3:A5").WrapText = $true
$ws.Range("C4:C5").WrapText = $true

# Row 3's Description cell (C3) gets its own distinct look: left/top aligned,
# wrapped, shaded with the light "Background 1" theme fill.
$ws.Range("C3").WrapText = $true
$ws.Range("C3").HorizontalAlignment = -4131
$ws.Range("C3").VerticalAlignment = -4108
$ws.Range("C3").Interior.ThemeColor = 2

# Bold the "DRA landing screen" run inside C3's rich text.
$boldStart = $richPrefix.Length + 1
$boldLen = $richBold.Length
$ws.Range("C3").Characters($boldStart, $boldLen).Font.Bold = $true

# Row heights matching the new content.
$ws.Rows.Item(3).RowHeight = 34.5
$ws.Rows.Item(4).RowHeight = 120
$ws.Rows.Item(5).RowHeight = 45

# Selection moves to D4 after the edit (as recorded in the workbook view).
$ws.Range("D4").Select()

Write-Host "done"
